# Auto-generated PowerShell/COM script to add the 'Query' column
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Order to Cash")
$ws2 = $wb.Worksheets.Item("Inventory")
$ws3 = $wb.Worksheets.Item("Finance")
$ws4 = $wb.Worksheets.Item("Procurement")

# ---- Order to Cash (sheet1): new column J ----
$ws1.Range("A1").Copy() | Out-Null
$ws1.Range("J1").PasteSpecial(-4122) | Out-Null
$ws1.Range("J1").Value = "Query"
$ws1.Range("G2").Copy() | Out-Null
$ws1.Range("J2").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J3").PasteSpecial(-4122) | Out-Null
$ws1.Range("J3").Value = "〇"
$ws1.Range("G2").Copy() | Out-Null
$ws1.Range("J4").PasteSpecial(-4122) | Out-Null
$ws1.Range("J4").Value = "〇"
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J5").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J6").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J7").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J8").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J9").PasteSpecial(-4122) | Out-Null
$ws1.Range("J9").Value = "〇"
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J10").PasteSpecial(-4122) | Out-Null
$ws1.Range("J10").Value = "〇"
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J11").PasteSpecial(-4122) | Out-Null
$ws1.Range("J11").Value = "〇"
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J12").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J13").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J14").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J15").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J16").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J17").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J18").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J19").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J20").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J21").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J22").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J23").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J24").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J25").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J26").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J27").PasteSpecial(-4122) | Out-Null
$ws1.Range("J27").Value = "〇"
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J28").PasteSpecial(-4122) | Out-Null
$ws1.Range("J28").Value = "〇"
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J29").PasteSpecial(-4122) | Out-Null
$ws1.Range("J29").Value = "〇"
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J30").PasteSpecial(-4122) | Out-Null
$ws1.Range("J30").Value = "〇"
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J31").PasteSpecial(-4122) | Out-Null
$ws1.Range("J31").Value = "〇"
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J32").PasteSpecial(-4122) | Out-Null
$ws1.Range("J32").Value = "〇"
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J33").PasteSpecial(-4122) | Out-Null
$ws1.Range("B5").Copy() | Out-Null
$ws1.Range("J34").PasteSpecial(-4122) | Out-Null
$ws1.Range("G35").Copy() | Out-Null
$ws1.Range("J35").PasteSpecial(-4122) | Out-Null

# ---- Inventory (sheet2): new column J ----
$ws3.Range("J1").Copy() | Out-Null
$ws2.Range("J1").PasteSpecial(-4122) | Out-Null
$ws2.Range("J1").Value = "Query"
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J2").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J3").PasteSpecial(-4122) | Out-Null
$ws2.Range("J3").Value = "〇"
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J4").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J5").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J6").PasteSpecial(-4122) | Out-Null
$ws2.Range("J6").Value = "〇"
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J7").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J8").PasteSpecial(-4122) | Out-Null
$ws2.Range("J8").Value = "〇"
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J9").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J10").PasteSpecial(-4122) | Out-Null
$ws2.Range("J10").Value = "〇"
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J11").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J12").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J13").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J14").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J15").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J16").PasteSpecial(-4122) | Out-Null
$ws2.Range("J16").Value = "〇"
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J17").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J18").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J19").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J20").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J21").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J22").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J23").PasteSpecial(-4122) | Out-Null
$ws2.Range("J23").Value = "〇"
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J24").PasteSpecial(-4122) | Out-Null
$ws2.Range("J24").Value = "〇"
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J25").PasteSpecial(-4122) | Out-Null
$ws2.Range("J25").Value = "〇"
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J26").PasteSpecial(-4122) | Out-Null
$ws2.Range("J26").Value = "〇"
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J27").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J28").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J29").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J30").PasteSpecial(-4122) | Out-Null
$ws2.Range("J30").Value = "〇"
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("J31").PasteSpecial(-4122) | Out-Null

# ---- Finance (sheet3): new column K ----
$ws3.Range("J1").Copy() | Out-Null
$ws3.Range("K1").PasteSpecial(-4122) | Out-Null
$ws3.Range("K1").Value = "Query"
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K2").PasteSpecial(-4122) | Out-Null
$ws3.Range("K2").Value = "〇"
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K3").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K4").PasteSpecial(-4122) | Out-Null
$ws3.Range("K4").Value = "〇"
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K5").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K6").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K7").PasteSpecial(-4122) | Out-Null
$ws3.Range("K7").Value = "〇"
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K8").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K9").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K10").PasteSpecial(-4122) | Out-Null
$ws3.Range("K10").Value = "〇"
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K11").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K12").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K13").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K14").PasteSpecial(-4122) | Out-Null
$ws3.Range("K14").Value = "〇"
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K15").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K16").PasteSpecial(-4122) | Out-Null
$ws3.Range("K16").Value = "〇"
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K17").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K18").PasteSpecial(-4122) | Out-Null
$ws3.Range("K18").Value = "〇"
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K19").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K20").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K21").PasteSpecial(-4122) | Out-Null
$ws3.Range("K21").Value = "〇"
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K22").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K23").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K24").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K25").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K26").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K27").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K28").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K29").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K30").PasteSpecial(-4122) | Out-Null
$ws3.Range("K30").Value = "〇"
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K31").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K32").PasteSpecial(-4122) | Out-Null
$ws3.Range("K32").Value = "〇"
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K33").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K34").PasteSpecial(-4122) | Out-Null
$ws3.Range("K34").Value = "〇"
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K35").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K36").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K37").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K38").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K39").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K40").PasteSpecial(-4122) | Out-Null
$ws3.Range("G2").Copy() | Out-Null
$ws3.Range("K41").PasteSpecial(-4122) | Out-Null
$ws3.Range("G42").Copy() | Out-Null
$ws3.Range("K42").PasteSpecial(-4122) | Out-Null

# ---- Procurement (sheet4): new column J ----
$ws3.Range("J1").Copy() | Out-Null
$ws4.Range("J1").PasteSpecial(-4122) | Out-Null
$ws4.Range("J1").Value = "Query"
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J2").PasteSpecial(-4122) | Out-Null
$ws4.Range("J2").Value = "〇"
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J3").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J4").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J5").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J6").PasteSpecial(-4122) | Out-Null
$ws4.Range("J6").Value = "〇"
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J7").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J8").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J9").PasteSpecial(-4122) | Out-Null
$ws4.Range("J9").Value = "〇"
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J10").PasteSpecial(-4122) | Out-Null
$ws4.Range("J10").Value = "〇"
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J11").PasteSpecial(-4122) | Out-Null
$ws4.Range("J11").Value = "〇"
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J12").PasteSpecial(-4122) | Out-Null
$ws4.Range("J12").Value = "〇"
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J13").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J14").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J15").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J16").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J17").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J18").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J19").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J20").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J21").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J22").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J23").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J24").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J25").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J26").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J27").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J28").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J29").PasteSpecial(-4122) | Out-Null
$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("J30").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---- AutoFilter ranges ----
$ws1.AutoFilterMode = $false
$ws1.Range("A1:J34").AutoFilter() | Out-Null
$ws2.AutoFilterMode = $false
$ws2.Range("A1:J31").AutoFilter() | Out-Null
$ws3.AutoFilterMode = $false
$ws3.Range("A1:K41").AutoFilter() | Out-Null
$ws4.AutoFilterMode = $false
$ws4.Range("A1:J30").AutoFilter() | Out-Null

# ---- Defined names (_FilterDatabase) ----
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Order to Cash!_FilterDatabase") { $n.RefersTo = "='Order to Cash'!`$A`$1:`$J`$34" }
    if ($n.Name -eq "Inventory!_FilterDatabase") { $n.RefersTo = "=Inventory!`$A`$1:`$J`$31" }
    if ($n.Name -eq "Finance!_FilterDatabase") { $n.RefersTo = "=Finance!`$A`$1:`$K`$41" }
    if ($n.Name -eq "Procurement!_FilterDatabase") { $n.RefersTo = "=Procurement!`$A`$1:`$J`$30" }
}

# ---- Active sheet / selection ----
$ws1.Activate()
$ws1.Select()
$ws1.Range("C25").Select()